# Update the "取得日時" (retrieved at) timestamp in column A for rows 2-11
# on the "ランサーズ" sheet from the old run timestamp to the new one,
# reflecting the latest append: 2025-12-02 18:38 JST.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-12-02 18:29:43"
$newTimestamp = "2025-12-02 18:38:21"

for ($row = 2; $row -le 11; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
